$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.731.66'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '1.742.89'
$ws.Range('E3').Value = '  -2.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3870'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3353'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.34'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.097'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07131'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.069'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.44%  '
$ws.Range('D15').Value = '1.741.81'
$ws.Range('E15').Value = '  -2.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.925'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001044'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06594'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.151'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.55%  '
$ws.Range('D23').Value = '27.733.60'
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.396'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.269'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.74%  '
$ws.Range('D29').Value = '1.937.54'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.265'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.039'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.733'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08656'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.91'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.92%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.073'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.96%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02237'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.73%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06023'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.01%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6370'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.24%  '
$ws.Range('B40').Value = 'WEMIXTOKEN'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.489'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2078'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.186'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.69%  '
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.849'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.802'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5887'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.959'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06920'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.139'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.82%  '
